$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 14999.75
$ws.Range("I51").Value = 14999.75
$ws.Range("K51").Value = 14999.75
$ws.Range("M51").Value = -14515.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4289.837
$ws.Range("I32").Value = 3195.8936
$ws.Range("J32").Value = 29997.5
$ws.Range("K32").Value = 3195.8936
$ws.Range("L32").Value = 29997.5
$ws.Range("M32").Value = -2908.8936
$ws.Range("N32").Value = -30571.5
$ws.Range("H61").Value = 1851.3572
$ws.Range("I61").Value = 1729.1666
$ws.Range("K61").Value = 1729.1666
$ws.Range("M61").Value = -1517.1666
$ws.Range("H82").Value = 80181
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 80181
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 80181
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -80903
$ws.Range("H85").Value = 80181
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 80181
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 80181
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -82677
$ws.Range("H110").Value = 1973
$ws.Range("I110").Value = 1709.625
$ws.Range("K110").Value = 1709.625
$ws.Range("M110").Value = 335.375
$ws.Range("H122").Value = 1826.25
$ws.Range("I122").Value = 1826.25
$ws.Range("K122").Value = 5478.75
$ws.Range("M122").Value = -3028.75
$ws.Range("H132").Value = 3608.389
$ws.Range("I132").Value = 2699.111
$ws.Range("J132").Value = 4517.6665
$ws.Range("K132").Value = 8097.333
$ws.Range("L132").Value = 13552.9995
$ws.Range("M132").Value = -5567.333
$ws.Range("N132").Value = -18612.9995
$ws.Range("H136").Value = 1851.3572
$ws.Range("I136").Value = 1729.1666
$ws.Range("K136").Value = 5187.4998
$ws.Range("M136").Value = -2637.4998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5411.5884
$ws.Range("J20").Value = 8424.857
$ws.Range("L20").Value = 8424.857
$ws.Range("N20").Value = -8918.857
$ws.Range("H94").Value = 1393.65
$ws.Range("I94").Value = 635.875
$ws.Range("J94").Value = 4424.75
$ws.Range("K94").Value = 635.875
$ws.Range("L94").Value = 4424.75
$ws.Range("M94").Value = -184.875
$ws.Range("N94").Value = -5326.75
$ws.Range("H134").Value = 2742.8333
$ws.Range("I134").Value = 2697.7144
$ws.Range("J134").Value = 2806
$ws.Range("K134").Value = 8093.1432
$ws.Range("L134").Value = 8418
$ws.Range("M134").Value = -5558.1432
$ws.Range("N134").Value = -13488

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 6014.6
$ws.Range("I69").Value = 2794
$ws.Range("J69").Value = 35000
$ws.Range("K69").Value = 2794
$ws.Range("L69").Value = 35000
$ws.Range("M69").Value = -2045
$ws.Range("N69").Value = -36498
$ws.Range("H72").Value = 6014.6
$ws.Range("I72").Value = 2794
$ws.Range("J72").Value = 35000
$ws.Range("K72").Value = 8382
$ws.Range("L72").Value = 105000
$ws.Range("M72").Value = -4638
$ws.Range("N72").Value = -112488

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 109.6
$ws.Range("I23").Value = 92
$ws.Range("J23").Value = 121.333336
$ws.Range("K23").Value = 276
$ws.Range("L23").Value = 364.000008
$ws.Range("M23").Value = -41
$ws.Range("N23").Value = -834.000008
$ws.Range("H33").Value = 211.57143
$ws.Range("I33").Value = 233
$ws.Range("K33").Value = 1398
$ws.Range("M33").Value = -1115
$ws.Range("H81").Value = 4006.5
$ws.Range("I81").Value = 3013
$ws.Range("K81").Value = 9039
$ws.Range("M81").Value = -7916
$ws.Range("H84").Value = 4006.5
$ws.Range("I84").Value = 3013
$ws.Range("K84").Value = 27117
$ws.Range("M84").Value = -21501
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H131").Value = 1714.5416
$ws.Range("J131").Value = 1944.5
$ws.Range("L131").Value = 5833.5
$ws.Range("N131").Value = -15913.5
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("M140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 127.052635
$ws.Range("I2").Value = 18.583334
$ws.Range("J2").Value = 313
$ws.Range("K2").Value = 18.583334
$ws.Range("L2").Value = 313
$ws.Range("M2").Value = 94.41666599999999
$ws.Range("N2").Value = -539
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 2000
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -1832
$ws.Range("N17").ClearContents()
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H70").Value = 7713014
$ws.Range("J70").Value = 21173.5
$ws.Range("L70").Value = 21173.5
$ws.Range("N70").Value = -21713.5
$ws.Range("H73").Value = 7713014
$ws.Range("J73").Value = 21173.5
$ws.Range("L73").Value = 21173.5
$ws.Range("N73").Value = -23045.5
$ws.Range("H93").Value = 75044.45
$ws.Range("J93").Value = 75044.45
$ws.Range("L93").Value = 75044.45
$ws.Range("N93").Value = -78788.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3275.6155
$ws.Range("J46").Value = 3409.2778
$ws.Range("L46").Value = 3409.2778
$ws.Range("N46").Value = -3785.2778
$ws.Range("H55").Value = 627.7895
$ws.Range("I55").Value = 233.4
$ws.Range("J55").Value = 1066
$ws.Range("K55").Value = 233.4
$ws.Range("L55").Value = 1066
$ws.Range("M55").Value = -60.40000000000001
$ws.Range("N55").Value = -1412
$ws.Range("H74").Value = 32549
$ws.Range("J74").Value = 39999.668
$ws.Range("L74").Value = 39999.668
$ws.Range("N74").Value = -41995.668
$ws.Range("H77").Value = 32549
$ws.Range("J77").Value = 39999.668
$ws.Range("L77").Value = 119999.004
$ws.Range("N77").Value = -129983.004
$ws.Range("H122").Value = 3200
$ws.Range("I122").Value = 3200
$ws.Range("K122").Value = 9600
$ws.Range("M122").Value = -7150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 6000
$ws.Range("I52").Value = 6000
$ws.Range("K52").Value = 6000
$ws.Range("M52").Value = -5774
$ws.Range("H75").Value = 24999.666
$ws.Range("I75").Value = 24999
$ws.Range("K75").Value = 24999
$ws.Range("M75").Value = -24063
$ws.Range("H78").Value = 24999.666
$ws.Range("I78").Value = 24999
$ws.Range("K78").Value = 74997
$ws.Range("M78").Value = -70317
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H107").Value = 594.9286
$ws.Range("I107").Value = 543.7143
$ws.Range("J107").Value = 646.1429000000001
$ws.Range("K107").Value = 1631.1429
$ws.Range("L107").Value = 1938.4287
$ws.Range("M107").Value = 288.8571000000002
$ws.Range("N107").Value = -5778.4287
